# Add 2022-Q3 data
#
# 1. A brand-new "2022-Q3" sheet is inserted right after "总计" (i.e. right
#    before the sheet that is currently named "2022-Q2"). It is built by
#    copying the existing "2022-Q2" sheet (so it inherits the same
#    formatting/styles/column layout) and then overwriting its data with the
#    new quarter's fund-holding figures.
# 2. The "总计" (summary) sheet gets a new row appended (2020-Q4 numbers,
#    which used to be the last row) and every existing row's date label /
#    numbers shift down by one quarter to make room for the new 2022-Q3
#    figures on row 2.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet by copying "2022-Q2", inserted before it
# ---------------------------------------------------------------------
$templateSheet = $wb.Worksheets.Item("2022-Q2")
$templateSheet.Copy($templateSheet)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# Row 2: 004497 / 前海开源多元策略灵活配置混合C
$q3.Cells.Item(2, 2).Value = "'004497"
$q3.Cells.Item(2, 3).Value = "前海开源多元策略灵活配置混合C"
$q3.Cells.Item(2, 4).Value = "'1.68"
$q3.Cells.Item(2, 5).Value = "'93.04"
$q3.Cells.Item(2, 6).Value = "'4.19"
$q3.Cells.Item(2, 7).Value = "'0.0704"
$q3.Cells.Item(2, 8).Value = 8

# Row 3: 160135 / 南方中证高铁产业指数（LOF）
$q3.Cells.Item(3, 2).Value = "'160135"
$q3.Cells.Item(3, 3).Value = "南方中证高铁产业指数（LOF）"
$q3.Cells.Item(3, 4).Value = "'1.84"
$q3.Cells.Item(3, 5).Value = "'95.01"
$q3.Cells.Item(3, 6).Value = "'2.41"
$q3.Cells.Item(3, 7).Value = "'0.0443"
$q3.Cells.Item(3, 8).Value = 10

# Row 4: 004496 / 前海开源多元策略灵活配置混合A
$q3.Cells.Item(4, 2).Value = "'004496"
$q3.Cells.Item(4, 3).Value = "前海开源多元策略灵活配置混合A"
$q3.Cells.Item(4, 4).Value = "'0.91"
$q3.Cells.Item(4, 5).Value = "'93.04"
$q3.Cells.Item(4, 6).Value = "'4.19"
$q3.Cells.Item(4, 7).Value = "'0.0381"
$q3.Cells.Item(4, 8).Value = 8

# Row 5: 160639 / 鹏华中证高铁产业指数（LOF）A
$q3.Cells.Item(5, 2).Value = "'160639"
$q3.Cells.Item(5, 3).Value = "鹏华中证高铁产业指数（LOF）A"
$q3.Cells.Item(5, 4).Value = "'0.75"
$q3.Cells.Item(5, 5).Value = "'94.62"
$q3.Cells.Item(5, 6).Value = "'2.39"
$q3.Cells.Item(5, 7).Value = "'0.0179"
$q3.Cells.Item(5, 8).Value = 10

# Row 6: 015678 / 鹏华中证高铁产业指数（LOF）C
$q3.Cells.Item(6, 2).Value = "'015678"
$q3.Cells.Item(6, 3).Value = "鹏华中证高铁产业指数（LOF）C"
$q3.Cells.Item(6, 4).Value = "'0.06"
$q3.Cells.Item(6, 5).Value = "'94.62"
$q3.Cells.Item(6, 6).Value = "'2.39"
$q3.Cells.Item(6, 7).Value = "'0.0014"
$q3.Cells.Item(6, 8).Value = 10

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: shift rows down one quarter and add
#    a new trailing row for 2020-Q4.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Extend the existing row 8 formatting/style down into the brand-new row 9
# by copying A8:D8 -> A9:D9, then overwrite the values below.
$total.Range("A8:D8").Copy($total.Range("A9:D9"))

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 5
$total.Cells.Item(2, 4).Value = 0.17

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(3, 2).Value = "2022-Q2"
$total.Cells.Item(3, 3).Value = 5
$total.Cells.Item(3, 4).Value = 0.29

$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(4, 2).Value = "2022-Q1"
$total.Cells.Item(4, 3).Value = 4
$total.Cells.Item(4, 4).Value = 0.23

$total.Cells.Item(5, 1).Value = 3
$total.Cells.Item(5, 2).Value = "2021-Q4"
$total.Cells.Item(5, 3).Value = 6
$total.Cells.Item(5, 4).Value = 0.31

$total.Cells.Item(6, 1).Value = 4
$total.Cells.Item(6, 2).Value = "2021-Q3"
$total.Cells.Item(6, 3).Value = 5
$total.Cells.Item(6, 4).Value = 0.11

$total.Cells.Item(7, 1).Value = 5
$total.Cells.Item(7, 2).Value = "2021-Q2"
$total.Cells.Item(7, 3).Value = 4
$total.Cells.Item(7, 4).Value = 0.26

$total.Cells.Item(8, 1).Value = 6
$total.Cells.Item(8, 2).Value = "2021-Q1"
$total.Cells.Item(8, 3).Value = 3
$total.Cells.Item(8, 4).Value = 0.2

$total.Cells.Item(9, 1).Value = 7
$total.Cells.Item(9, 2).Value = "2020-Q4"
$total.Cells.Item(9, 3).Value = 3
$total.Cells.Item(9, 4).Value = 0.19

# Restore the originally-active tab (last sheet, "2020-Q4") so the active
# sheet / tabSelected marker stays where it was before this edit instead of
# lingering on the newly created "2022-Q3" sheet.
$wb.Worksheets.Item("2020-Q4").Activate()
